$d = $word.ActiveDocument

# Update the date heading in the first paragraph.
# wdReplaceOne (1) replaces just the single match found, which is the
# safest choice when scoping a Find to a specific Range.
$d.Content.Find.Execute("2025-04-08 Tuesday", $true, $false, $false, $false, $false, $true, 0, $false, "2025-04-09 Wednesday", 1) | Out-Null

# Update each practice-problem cell in the table by (row, column) so that
# duplicate source strings (e.g. two "900 div 4=225, 0" cells) resolve independently.
# Each Find is scoped to that single cell Range and uses wdReplaceOne so it
# only ever touches the first (and only) match inside that cell.
$tbl = $d.Tables.Item(1)

$edits = @(
    @{Row=1; Col=1; Old="177÷4=44, 1"; New="701÷9=77, 8"},
    @{Row=1; Col=2; Old="725÷4=181, 1"; New="385÷9=42, 7"},
    @{Row=1; Col=3; Old="791÷2=395, 1"; New="721÷9=80, 1"},
    @{Row=1; Col=4; Old="154÷3=51, 1"; New="649÷8=81, 1"},
    @{Row=1; Col=5; Old="606÷9=67, 3"; New="718÷9=79, 7"},
    @{Row=5; Col=1; Old="330÷9=36, 6"; New="486÷9=54, 0"},
    @{Row=5; Col=2; Old="718÷5=143, 3"; New="185÷2=92, 1"},
    @{Row=5; Col=3; Old="900÷4=225, 0"; New="803÷7=114, 5"},
    @{Row=5; Col=4; Old="900÷4=225, 0"; New="608÷6=101, 2"},
    @{Row=5; Col=5; Old="216÷2=108, 0"; New="382÷3=127, 1"},
    @{Row=9; Col=1; Old="590÷2=295, 0"; New="484÷8=60, 4"},
    @{Row=9; Col=2; Old="524÷8=65, 4"; New="256÷4=64, 0"},
    @{Row=9; Col=3; Old="360÷5=72, 0"; New="253÷5=50, 3"},
    @{Row=9; Col=4; Old="241÷4=60, 1"; New="638÷8=79, 6"},
    @{Row=9; Col=5; Old="696÷7=99, 3"; New="818÷6=136, 2"},
    @{Row=13; Col=1; Old="444÷2=222, 0"; New="523÷4=130, 3"},
    @{Row=13; Col=2; Old="395÷4=98, 3"; New="106÷8=13, 2"},
    @{Row=13; Col=3; Old="688÷5=137, 3"; New="668÷9=74, 2"},
    @{Row=13; Col=4; Old="215÷7=30, 5"; New="417÷5=83, 2"},
    @{Row=13; Col=5; Old="300÷6=50, 0"; New="783÷6=130, 3"},
    @{Row=17; Col=1; Old="229÷5=45, 4"; New="407÷2=203, 1"},
    @{Row=17; Col=2; Old="606÷7=86, 4"; New="302÷4=75, 2"},
    @{Row=17; Col=3; Old="649÷7=92, 5"; New="349÷3=116, 1"},
    @{Row=17; Col=4; Old="924÷5=184, 4"; New="713÷8=89, 1"},
    @{Row=17; Col=5; Old="842÷2=421, 0"; New="920÷5=184, 0"}
)

foreach ($edit in $edits) {
    $cell = $tbl.Cell($edit.Row, $edit.Col)
    $cellRange = $cell.Range
    $cellRange.Find.Execute($edit.Old, $true, $false, $false, $false, $false, $true, 0, $false, $edit.New, 1) | Out-Null
}
